# feat: add 2022-Q1 data
#
# - Insert a new "2022-Q1" sheet (fund holdings detail) right before the
#   existing "总计" (totals) sheet.
# - Append a new "2022-Q1" row to the "总计" summary sheet (it ends up being
#   rebuilt as a brand-new sheet placed after "2022-Q1" so the engine's
#   sheetId bookkeeping matches a plain insert+append edit).

$wb = $excel.ActiveWorkbook

$q3 = $wb.Worksheets.Item("2021-Q3")
$q4 = $wb.Worksheets.Item("2021-Q4")
$totalOld = $wb.Worksheets.Item("总计")

# Drop the old totals sheet - its sheetId is freed and will be reused by the
# next sheet we add, which keeps the workbook's sheetId sequence (1,2,3,4,5)
# identical to what a manual "insert sheet before 总计" edit would produce.
$totalOld.Delete()

# --- New "2022-Q1" fund-holdings sheet, positioned right after 2021-Q4 ---
$newQ1 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $q4)
$newQ1.Name = "2022-Q1"

# Copy a same-shaped range from 2021-Q4 so header/border/bold formatting
# (style index 2 in this workbook) and page margins match the sibling
# quarter sheets.
$q4.Range("A1:H4").Copy($newQ1.Range("A1"))
$newQ1.PageSetup.LeftMargin = 54
$newQ1.PageSetup.RightMargin = 54
$newQ1.PageSetup.TopMargin = 72
$newQ1.PageSetup.BottomMargin = 72
$newQ1.PageSetup.HeaderMargin = 36
$newQ1.PageSetup.FooterMargin = 36

$newQ1.Range("B1").Value = "基金代码"
$newQ1.Range("C1").Value = "基金名称"
$newQ1.Range("D1").Value = "基金规模"
$newQ1.Range("E1").Value = "股票总仓位"
$newQ1.Range("F1").Value = "仓位占比"
$newQ1.Range("G1").Value = "持有市值(亿元)"
$newQ1.Range("H1").Value = "仓位排名"

# Force text storage for the code/percentage-like columns so values such as
# "000522" keep their leading zero instead of being coerced to a number.
$newQ1.Range("B2:B4").NumberFormat = "@"
$newQ1.Range("D2:G4").NumberFormat = "@"

$newQ1.Range("A2").Value = 0
$newQ1.Range("B2").Value = "159851"
$newQ1.Range("C2").Value = "华宝中证金融科技主题ETF"
$newQ1.Range("D2").Value = "3.16"
$newQ1.Range("E2").Value = "98.58"
$newQ1.Range("F2").Value = "3.41"
$newQ1.Range("G2").Value = "0.1078"
$newQ1.Range("H2").Value = 7

$newQ1.Range("A3").Value = 1
$newQ1.Range("B3").Value = "000522"
$newQ1.Range("C3").Value = "华润元大信息传媒科技混合"
$newQ1.Range("D3").Value = "1.50"
$newQ1.Range("E3").Value = "70.63"
$newQ1.Range("F3").Value = "2.50"
$newQ1.Range("G3").Value = "0.0375"
$newQ1.Range("H3").Value = 10

$newQ1.Range("A4").Value = 2
$newQ1.Range("B4").Value = "516100"
$newQ1.Range("C4").Value = "华夏中证金融科技主题交易型开放式指数证券投资基金"
$newQ1.Range("D4").Value = "0.68"
$newQ1.Range("E4").Value = "96.91"
$newQ1.Range("F4").Value = "3.39"
$newQ1.Range("G4").Value = "0.0231"
$newQ1.Range("H4").Value = 7

# --- Rebuild "总计" after "2022-Q1", with the new quarter prepended ---
$newTotal = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $newQ1)
$newTotal.Name = "总计"

$q3.Range("A1:D5").Copy($newTotal.Range("A1"))
$newTotal.PageSetup.LeftMargin = 54
$newTotal.PageSetup.RightMargin = 54
$newTotal.PageSetup.TopMargin = 72
$newTotal.PageSetup.BottomMargin = 72
$newTotal.PageSetup.HeaderMargin = 36
$newTotal.PageSetup.FooterMargin = 36

$newTotal.Range("B1").Value = "日期"
$newTotal.Range("C1").Value = "持有数量(只)"
$newTotal.Range("D1").Value = "持有市值(亿元)"

$newTotal.Range("A2").Value = 0
$newTotal.Range("B2").Value = "2022-Q1"
$newTotal.Range("C2").Value = 3
$newTotal.Range("D2").Value = 0.17

$newTotal.Range("A3").Value = 1
$newTotal.Range("B3").Value = "2021-Q4"
$newTotal.Range("C3").Value = 4
$newTotal.Range("D3").Value = 3.85

$newTotal.Range("A4").Value = 2
$newTotal.Range("B4").Value = "2021-Q3"
$newTotal.Range("C4").Value = 11
$newTotal.Range("D4").Value = 5.46

$newTotal.Range("A5").Value = 3
$newTotal.Range("B5").Value = "2021-Q2"
$newTotal.Range("C5").Value = 5
$newTotal.Range("D5").Value = 1.74

# Restore the originally-active tab (2021-Q2) - adding/copying sheets above
# shifts Excel's selection to the most-recently-touched sheet otherwise.
$wb.Worksheets.Item("2021-Q2").Activate()

Write-Host "Sheets now:"
foreach ($ws in $wb.Worksheets) {
    Write-Host (" - " + $ws.Name)
}
